$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 85, shifting existing rows 85-110 down to 86-111.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new weekly price observation.
$ws.Cells.Item(85, 1).Value = 6
$ws.Cells.Item(85, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(85, 3).Value = "Metropolitana"
$ws.Cells.Item(85, 4).Value = 44875
$ws.Cells.Item(85, 5).Value = 13
$ws.Cells.Item(85, 6).Value = 100114007
$ws.Cells.Item(85, 7).Value = "Jengibre"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 1400
$ws.Cells.Item(85, 11).Value = 13000
$ws.Cells.Item(85, 12).Value = 14000
$ws.Cells.Item(85, 13).Value = 13536
$ws.Cells.Item(85, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(85, 15).Value = "Perú"
$ws.Cells.Item(85, 16).Value = 1041
$ws.Cells.Item(85, 17).Value = 13
$ws.Cells.Item(85, 18).Value = "Hortaliza"

Write-Output "Inserted new row 85 and shifted subsequent rows down."
